$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '68.772.80'
$ws.Cells.Item(2, 5).Value = '  -2.27%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.459.46'
$ws.Cells.Item(3, 5).Value = '  -4.40%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.01%  '

# Row 5
$ws.Cells.Item(5, 4).Value = "'573.38"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -4.76%  '

# Row 6
$ws.Cells.Item(6, 4).Value = "'190.36"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -3.06%  '

# Row 7
$ws.Cells.Item(7, 2).Value = 'XRP'
$ws.Cells.Item(7, 3).Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Cells.Item(7, 4).Value = "'0.605"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value = '  -3.31%  '

# Row 8
$ws.Cells.Item(8, 2).Value = 'LidoStakedEther'
$ws.Cells.Item(8, 3).Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Cells.Item(8, 4).Value = '3.448.70'
$ws.Cells.Item(8, 5).Value = '  -4.36%  '

# Row 9
$ws.Cells.Item(9, 5).Value = '  +0.07%  '

# Row 10
$ws.Cells.Item(10, 4).Value = "'0.202"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  -4.40%  '

# Row 11
$ws.Cells.Item(11, 4).Value = "'0.614"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  -5.10%  '

# Row 12
$ws.Cells.Item(12, 4).Value = "'51.09"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -4.07%  '

# Row 13
$ws.Cells.Item(13, 4).Value = "'0.0000284"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -6.51%  '

# Row 14
$ws.Cells.Item(14, 4).Value = "'9.05"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -5.36%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.999.89'
$ws.Cells.Item(15, 5).Value = '  -4.73%  '

# Row 16
$ws.Cells.Item(16, 4).Value = "'637.30"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  +5.48%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '68.544.36'
$ws.Cells.Item(17, 5).Value = '  -2.74%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.461.12'
$ws.Cells.Item(18, 5).Value = '  -3.72%  '

# Row 19
$ws.Cells.Item(19, 4).Value = "'12.28"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -4.79%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -2.50%  '

# Row 21
$ws.Cells.Item(21, 4).Value = "'18.09"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -5.24%  '

# Row 22
$ws.Cells.Item(22, 4).Value = "'0.938"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -6.06%  '

# Row 23
$ws.Cells.Item(23, 4).Value = "'17.85"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -1.64%  '

# Row 24
$ws.Cells.Item(24, 4).Value = "'5.35"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +3.44%  '

# Row 25
$ws.Cells.Item(25, 4).Value = "'99.15"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -3.78%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  -7.69%  '

# Row 27
$ws.Cells.Item(27, 4).Value = "'2.83"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  -5.47%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  +1.94%  '

# Row 29
$ws.Cells.Item(29, 4).Value = "'9.79"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -8.14%  '

# Row 30
$ws.Cells.Item(30, 4).Value = "'9.18"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  -5.12%  '

# Row 31
$ws.Cells.Item(31, 4).Value = "'32.26"
$ws.Cells.Item(31, 4).Style = 'Normal'
$ws.Cells.Item(31, 5).Value = '  -4.53%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  -11.56%  '

# Row 33
$ws.Cells.Item(33, 4).Value = "'6.70"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -8.20%  '

# Row 34
$ws.Cells.Item(34, 4).Value = "'11.54"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  -6.00%  '

# Row 35
$ws.Cells.Item(35, 4).Value = "'60.85"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -4.12%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -7.56%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +0.06%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '3.663.85'
$ws.Cells.Item(38, 5).Value = '  -6.17%  '

# Row 39
$ws.Cells.Item(39, 2).Value = 'PEPE'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(39, 4).Value = '0.0₃0780'
$ws.Cells.Item(39, 5).Value = '  -11.54%  '

# Row 40
$ws.Cells.Item(40, 2).Value = 'Bittensor'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(40, 4).Value = "'502.49"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -3.01%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -0.79%  '

# Row 42
$ws.Cells.Item(42, 4).Value = "'2.89"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -6.14%  '

# Row 43
$ws.Cells.Item(43, 4).Value = "'0.366"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -5.57%  '

# Row 44
$ws.Cells.Item(44, 4).Value = "'34.27"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -6.99%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -2.03%  '

# Row 46
$ws.Cells.Item(46, 4).Value = "'3.45"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +65.04%  '

# Row 47
$ws.Cells.Item(47, 4).Value = "'0.0435"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -5.24%  '

# Row 48
$ws.Cells.Item(48, 5).Value = '  -3.77%  '

# Row 49
$ws.Cells.Item(49, 4).Value = "'2.80"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -3.42%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -4.45%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -0.42%  '
